$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-11-27 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-28 Friday", 2) | Out-Null

# Update each answer cell in the table, row by row, column by column
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "10+43=53"
$t.Cell(1, 2).Range.Text = "48+43=91"
$t.Cell(1, 3).Range.Text = "39-2=37"
$t.Cell(1, 4).Range.Text = "72-62=10"
$t.Cell(1, 5).Range.Text = "92-43=49"

$t.Cell(2, 1).Range.Text = "49+8=57"
$t.Cell(2, 2).Range.Text = "20+8=28"
$t.Cell(2, 3).Range.Text = "12+10=22"
$t.Cell(2, 4).Range.Text = "99-10=89"
$t.Cell(2, 5).Range.Text = "30+60=90"

$t.Cell(3, 1).Range.Text = "87-33=54"
$t.Cell(3, 2).Range.Text = "68-3=65"
$t.Cell(3, 3).Range.Text = "4+8=12"
$t.Cell(3, 4).Range.Text = "19+59=78"
$t.Cell(3, 5).Range.Text = "59+1=60"

$t.Cell(4, 1).Range.Text = "38+10=48"
$t.Cell(4, 2).Range.Text = "82-25=57"
$t.Cell(4, 3).Range.Text = "55-41=14"
$t.Cell(4, 4).Range.Text = "81-18=63"
$t.Cell(4, 5).Range.Text = "80-22=58"

$t.Cell(5, 1).Range.Text = "56-40=16"
$t.Cell(5, 2).Range.Text = "93-87=6"
$t.Cell(5, 3).Range.Text = "34-32=2"
$t.Cell(5, 4).Range.Text = "4+53=57"
$t.Cell(5, 5).Range.Text = "15+54=69"

$t.Cell(6, 1).Range.Text = "97-23=74"
$t.Cell(6, 2).Range.Text = "61-53=8"
$t.Cell(6, 3).Range.Text = "54-48=6"
$t.Cell(6, 4).Range.Text = "53-10=43"
$t.Cell(6, 5).Range.Text = "4+55=59"

$t.Cell(7, 1).Range.Text = "92-3=89"
$t.Cell(7, 2).Range.Text = "92-91=1"
$t.Cell(7, 3).Range.Text = "31+63=94"
$t.Cell(7, 4).Range.Text = "60+14=74"
$t.Cell(7, 5).Range.Text = "80+8=88"

$t.Cell(8, 1).Range.Text = "68+5=73"
$t.Cell(8, 2).Range.Text = "91-86=5"
$t.Cell(8, 3).Range.Text = "36+35=71"
$t.Cell(8, 4).Range.Text = "45+32=77"
$t.Cell(8, 5).Range.Text = "33+3=36"

$t.Cell(9, 1).Range.Text = "91-53=38"
$t.Cell(9, 2).Range.Text = "45+15=60"
$t.Cell(9, 3).Range.Text = "31+23=54"
$t.Cell(9, 4).Range.Text = "12+37=49"
$t.Cell(9, 5).Range.Text = "60-47=13"

$t.Cell(10, 1).Range.Text = "39-27=12"
$t.Cell(10, 2).Range.Text = "90-46=44"
$t.Cell(10, 3).Range.Text = "99-32=67"
$t.Cell(10, 4).Range.Text = "67-67=0"
$t.Cell(10, 5).Range.Text = "16+14=30"

$t.Cell(11, 1).Range.Text = "8+55=63"
$t.Cell(11, 2).Range.Text = "96-46=50"
$t.Cell(11, 3).Range.Text = "41-15=26"
$t.Cell(11, 4).Range.Text = "18+34=52"
$t.Cell(11, 5).Range.Text = "9+82=91"

$t.Cell(12, 1).Range.Text = "51-21=30"
$t.Cell(12, 2).Range.Text = "2+32=34"
$t.Cell(12, 3).Range.Text = "3+50=53"
$t.Cell(12, 4).Range.Text = "47-42=5"
$t.Cell(12, 5).Range.Text = "46+52=98"

$t.Cell(13, 1).Range.Text = "2+88=90"
$t.Cell(13, 2).Range.Text = "91-1=90"
$t.Cell(13, 3).Range.Text = "41+1=42"
$t.Cell(13, 4).Range.Text = "91-50=41"
$t.Cell(13, 5).Range.Text = "46-37=9"

$t.Cell(14, 1).Range.Text = "7+72=79"
$t.Cell(14, 2).Range.Text = "48-37=11"
$t.Cell(14, 3).Range.Text = "68+5=73"
$t.Cell(14, 4).Range.Text = "38-0=38"
$t.Cell(14, 5).Range.Text = "17+57=74"

$t.Cell(15, 1).Range.Text = "30-28=2"
$t.Cell(15, 2).Range.Text = "50+26=76"
$t.Cell(15, 3).Range.Text = "33+57=90"
$t.Cell(15, 4).Range.Text = "40+25=65"
$t.Cell(15, 5).Range.Text = "40+11=51"

$t.Cell(16, 1).Range.Text = "65-45=20"
$t.Cell(16, 2).Range.Text = "11+8=19"
$t.Cell(16, 3).Range.Text = "94-60=34"
$t.Cell(16, 4).Range.Text = "65+0=65"
$t.Cell(16, 5).Range.Text = "45-41=4"

$t.Cell(17, 1).Range.Text = "56+22=78"
$t.Cell(17, 2).Range.Text = "81-0=81"
$t.Cell(17, 3).Range.Text = "34-22=12"
$t.Cell(17, 4).Range.Text = "22+21=43"
$t.Cell(17, 5).Range.Text = "78-24=54"

$t.Cell(18, 1).Range.Text = "52-7=45"
$t.Cell(18, 2).Range.Text = "96-24=72"
$t.Cell(18, 3).Range.Text = "29+42=71"
$t.Cell(18, 4).Range.Text = "23+38=61"
$t.Cell(18, 5).Range.Text = "63-14=49"

$t.Cell(19, 1).Range.Text = "6-1=5"
$t.Cell(19, 2).Range.Text = "3+34=37"
$t.Cell(19, 3).Range.Text = "36-11=25"
$t.Cell(19, 4).Range.Text = "26-18=8"
$t.Cell(19, 5).Range.Text = "81-19=62"

$t.Cell(20, 1).Range.Text = "12+42=54"
$t.Cell(20, 2).Range.Text = "18+67=85"
$t.Cell(20, 3).Range.Text = "90-23=67"
$t.Cell(20, 4).Range.Text = "74-35=39"
$t.Cell(20, 5).Range.Text = "87-81=6"
